$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLogin")

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "damager"

$ws.Range("A4").Select()
